# Build site at 2023-04-12 14:53:07 UTC
# This script fixes the LOB1258 course-info sheet:
#  - realigns column A labels with their corresponding B/C content
#  - inserts two new rows for the two responsible-professor names
#  - adds the missing Portuguese objective / short-syllabus / syllabus / bibliography text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows right after row 12 ("Docentes responsaveis:") ---
# These will hold the two professor names (previously misplaced further down).
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(14).Insert()

# The insert operation copies formatting (incl. the bold "A" style) into the
# new rows; column A must stay empty on these two rows, so clear it.
$ws.Cells.Item(13,1).Clear()
$ws.Cells.Item(14,1).Clear()

# Give B13/C13 and B14/C14 the same look as the other content cells
# (style "2" for column B, style "3" for column C) by copying formats
# from an existing, correctly-styled row instead of recreating styles.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Fill in the values ---

# Objetivos: (row 10) previously held the professor name by mistake;
# now holds the actual Portuguese objective text.
$ws.Range("B10").Value2 = "Fornecer ao aluno os conceitos básicos de hidráulica aplicadas ao meio ambiente."
$ws.Range("C10").Value2 = "Fornecer ao aluno os conceitos básicos de hidráulica aplicadas ao meio ambiente."

# Docentes responsaveis: (new rows 13 and 14)
$ws.Range("B13").Value2 = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Range("C13").Value2 = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Range("B14").Value2 = "7455355 - Robson da Silva Rocha"
$ws.Range("C14").Value2 = "7455355 - Robson da Silva Rocha"

# Programa resumido: (row 15, shifted down from 13) gets the Portuguese short syllabus
$ws.Range("B15").Value2 = "Escoamento permanente uniforme em condutos forçados e perdas de carga nas tubulações. Instalações de bombeamento e bombas hidráulicas. Escoamento permanente uniforme em condutos livres, resistência ao escoamento e perdas de carga nos canais naturais e artificiais. Estudo da carga específica em canais."
$ws.Range("C15").Value2 = "Escoamento permanente uniforme em condutos forçados e perdas de carga nas tubulações. Instalações de bombeamento e bombas hidráulicas. Escoamento permanente uniforme em condutos livres, resistência ao escoamento e perdas de carga nos canais naturais e artificiais. Estudo da carga específica em canais."

# Programa: (row 17, shifted down from 15) gets the Portuguese syllabus
$ws.Range("B17").Value2 = "- Hidrostática,- piezometria,- conservação da massa e quantidade de movimento,- Escoamentos Permanentes em Condutos Forçados,- Resistência ao Escoamento e Perdas de Carga,- Bombas e sistemas de recalque,- Escoamento Permanente Uniforme em Condutos Livres,- Resistência ao escoamento e Perdas de Carga,- Canais regulares e naturais,- Carga Específica,- Escoamento Permanente Gradualmente Variado,- Cálculo da linha d’água,- Ressalto Hidráulico."
$ws.Range("C17").Value2 = "- Hidrostática,- piezometria,- conservação da massa e quantidade de movimento,- Escoamentos Permanentes em Condutos Forçados,- Resistência ao Escoamento e Perdas de Carga,- Bombas e sistemas de recalque,- Escoamento Permanente Uniforme em Condutos Livres,- Resistência ao escoamento e Perdas de Carga,- Canais regulares e naturais,- Carga Específica,- Escoamento Permanente Gradualmente Variado,- Cálculo da linha d’água,- Ressalto Hidráulico."

# Método: (row 20, shifted down from 18) gets the method text (shifted up from old row 19)
$ws.Range("B20").Value2 = "Aulas teóricas e práticas, trabalhos de campo e exercícios dirigidos.Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios."
$ws.Range("C20").Value2 = "Aulas teóricas e práticas, trabalhos de campo e exercícios dirigidos.Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios."

# Critério: (row 21) gets the weighted-average text (shifted up from old row 20)
$ws.Range("B21").Value2 = "Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios."
$ws.Range("C21").Value2 = "Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios."

# Norma de recuperação: (row 22) gets the recovery-exam text (shifted up from old row 21)
$ws.Range("B22").Value2 = "1 (uma) prova de recuperação (R), sendo considerado aprovado se 0,5(NF + R) >= 5,0."
$ws.Range("C22").Value2 = "1 (uma) prova de recuperação (R), sendo considerado aprovado se 0,5(NF + R) >= 5,0."

# Bibliografia: (row 23) gets the full bibliography text (new content)
$ws.Range("B23").Value2 = 'Hauer, B.F. Lamberti, G.A. – Methods in Stream Ecology, 896 p, Academic Press, ISBN 0123329078, 2006Alfredini, Paolo; "Obras e Gestão de Portos e Costas - A Técnica Aliada ao Enfoque Logístico e Ambiental".Editora Edgard Blucher, São Paulo, 2005.ALFREDINI, P. - “Obras e Gestão de Portos e Costas - A Técnica Aliada ao Enfoque Logístico e Ambiental”. Editora Edgard Blucher, São Paulo, 2005.LENCASTRE, A. – “Hidráulica Geral”. Edição do Autor, Lisboa, 1996.OPEN UNIVERSITY COURSE TEAM – “Waves, tides and shallow water processes”. 1998ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 1999. (Coordenação: Luiz Di Bernardo). ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 200 (Coordenação: Marco A.P. Reali).'
$ws.Range("C23").Value2 = 'Hauer, B.F. Lamberti, G.A. – Methods in Stream Ecology, 896 p, Academic Press, ISBN 0123329078, 2006Alfredini, Paolo; "Obras e Gestão de Portos e Costas - A Técnica Aliada ao Enfoque Logístico e Ambiental".Editora Edgard Blucher, São Paulo, 2005.ALFREDINI, P. - “Obras e Gestão de Portos e Costas - A Técnica Aliada ao Enfoque Logístico e Ambiental”. Editora Edgard Blucher, São Paulo, 2005.LENCASTRE, A. – “Hidráulica Geral”. Edição do Autor, Lisboa, 1996.OPEN UNIVERSITY COURSE TEAM – “Waves, tides and shallow water processes”. 1998ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 1999. (Coordenação: Luiz Di Bernardo). ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 200 (Coordenação: Marco A.P. Reali).'

Write-Host "Done"
